$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string text is interned in the order it is first written, and the
# target file expects "target sisa" (B5) before "trg tmbh" (H2), so write B5 first.
$ws.Range("B5").Value = "target sisa"

# Row 2 - new "trg tmbh" header
$ws.Range("H2").Value = "trg tmbh"

# Row 3 - updated target/day figures + new helper cells (J3:O3)
$ws.Range("C3").Value = 11
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 10
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 17
$ws.Range("M3").Value = 3
$ws.Range("N3").Value = 3
$ws.Range("O3").Value = 3

# Row 4 - updated figures; H3 (old leftover value) is fully removed; H4/J4 become formulas
$ws.Range("C4").Value = 3
$ws.Range("F4").Value = 4
$ws.Range("H3").Clear()
$ws.Range("H4").Formula = "=C5-F5-F4"
$ws.Range("J4").Formula = "=J3/K3"

# Row 5 - new "target sisa" row with its own formulas
$ws.Range("C5").Formula = "=C3-G6"
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 1
$ws.Range("H5").Formula = "=C5-F5"

# Row 6 - new running total of the checkout column
$ws.Range("G6").Formula = "=SUM(G3:G5)"

# Selection moves to G7 (and the saved top-left cell resets to the sheet default)
$ws.Range("G7").Select()
